$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 20 with forecast data, matching existing row formatting
$ws.Range("A19").Copy($ws.Range("A20"))
$ws.Range("A20").Value = 45986

$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 2.043309689777173
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 1.199077969291551
